$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (prices / 1h volume %) - GitHub Actions scheduled update

$ws.Range("D2").Value = '62.393.92'
$ws.Range("E2").Value = '  +2.67%  '
$ws.Range("D3").Value = '3.432.07'
$ws.Range("E3").Value = '  +1.88%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '407.32'
$ws.Range("E5").Value = '  +0.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.10'
$ws.Range("E6").Value = '  +4.92%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.596'
$ws.Range("E7").Value = '  -1.87%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +3.64%  '
$ws.Range("E10").Value = '  +10.00%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.96'
$ws.Range("E11").Value = '  +0.23%  '
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.90'
$ws.Range("E13").Value = '  +1.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.41'
$ws.Range("E14").Value = '  -0.82%  '
$ws.Range("D15").Value = '3.434.06'
$ws.Range("E15").Value = '  +1.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '11.67'
$ws.Range("E16").Value = '  +1.71%  '
$ws.Range("D17").Value = '62.205.88'
$ws.Range("E17").Value = '  +2.30%  '
$ws.Range("E18").Value = '  +0.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000151'
$ws.Range("E19").Value = '  +12.89%  '
$ws.Range("E20").Value = '  -1.63%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '84.31'
$ws.Range("E21").Value = '  +2.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '312.44'
$ws.Range("E22").Value = '  +2.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.79'
$ws.Range("E23").Value = '  -1.78%  '
$ws.Range("E24").Value = '  +0.81%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.75'
$ws.Range("E25").Value = '  +1.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '29.75'
$ws.Range("E26").Value = '  +1.06%  '
$ws.Range("E27").Value = '  -3.95%  '
$ws.Range("E28").Value = '  +5.42%  '
$ws.Range("E29").Value = '  +8.58%  '
$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.172'
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '43.81'
$ws.Range("E31").Value = '  +2.76%  '
$ws.Range("E32").Value = '  -0.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.35'
$ws.Range("E33").Value = '  -2.80%  '
$ws.Range("E34").Value = '  +0.09%  '
$ws.Range("E35").Value = '  +0.97%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '51.68'
$ws.Range("E36").Value = '  -0.85%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.997'
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("E38").Value = '  +1.91%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.316'
$ws.Range("E40").Value = '  +11.80%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '141.86'
$ws.Range("E41").Value = '  +4.59%  '
$ws.Range("E42").Value = '  +0.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.97'
$ws.Range("E43").Value = '  -1.53%  '
$ws.Range("E44").Value = '  +0.43%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.78'
$ws.Range("E45").Value = '  -0.32%  '
$ws.Range("E46").Value = '  +0.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '21.36'
$ws.Range("E47").Value = '  -2.04%  '
$ws.Range("D48").Value = '2.102.95'
$ws.Range("E48").Value = '  -1.54%  '
$ws.Range("E49").Value = '  -1.17%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.96'
$ws.Range("E50").Value = '  +3.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.72'
$ws.Range("E51").Value = '  +18.64%  '
